$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Generated on" timestamp text in A2
$ws.Range("A2").Value = "Generated on: Fri Mar 26 23:51:27 WEST 2021"

# Refresh the randomly-generated "Date Time" values in column G
$ws.Range("G5").Value = 33026.48819444444
$ws.Range("G6").Value = 42047.66805555556

# Drop the now-unused "Percent" column values in H, clearing both the
# content and the formatting so the cells become fully empty
$ws.Range("H5:H6").ClearContents()
$ws.Range("H5:H6").ClearFormats()
